$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044238316943338
$ws.Range("D2").Value = 1.049089240646592
$ws.Range("E2").Value = 1.051723471404605
$ws.Range("F2").Value = 1.061557240437975
$ws.Range("I2").Value = 1.038991381183959
$ws.Range("J2").Value = 1.049304472907919
$ws.Range("K2").Value = 1.051847270569546
$ws.Range("L2").Value = 1.054474182563662
$ws.Range("M2").Value = 1.064280977605849
$ws.Range("N2").Value = 1.020258757610625
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045415301784194
$ws.Range("D3").Value = 1.049988191798903
$ws.Range("E3").Value = 1.052755896734827
$ws.Range("F3").Value = 1.062644639907438
$ws.Range("I3").Value = 1.039241597475147
$ws.Range("J3").Value = 1.050127505318546
$ws.Range("K3").Value = 1.052558173561454
$ws.Range("L3").Value = 1.055318746365217
$ws.Range("M3").Value = 1.065182334213784
$ws.Range("N3").Value = 1.020536017625811
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046176727330375
$ws.Range("D4").Value = 1.050569403113178
$ws.Range("E4").Value = 1.053424117441139
$ws.Range("F4").Value = 1.063348362446693
$ws.Range("I4").Value = 1.039401571391594
$ws.Range("J4").Value = 1.050659369771291
$ws.Range("K4").Value = 1.053017073169354
$ws.Range("L4").Value = 1.055864802486105
$ws.Range("M4").Value = 1.065765078046206
$ws.Range("N4").Value = 1.020715078685926
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046496792985091
$ws.Range("D5").Value = 1.050813631873336
$ws.Range("E5").Value = 1.053705079070315
$ws.Range("F5").Value = 1.06364423226291
$ws.Range("I5").Value = 1.039468361664803
$ws.Range("J5").Value = 1.050882800239402
$ws.Range("K5").Value = 1.053209730895426
$ws.Range("L5").Value = 1.056094261149412
$ws.Range("M5").Value = 1.06600994555994
$ws.Range("N5").Value = 1.020790273512054
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046550531244828
$ws.Range("D6").Value = 1.050854632363591
$ws.Range("E6").Value = 1.053752256216327
$ws.Range("F6").Value = 1.063693911574951
$ws.Range("I6").Value = 1.039479548907299
$ws.Range("J6").Value = 1.050920305515065
$ws.Range("K6").Value = 1.053242063519119
$ws.Range("L6").Value = 1.056132782219799
$ws.Range("M6").Value = 1.066051052982697
$ws.Range("N6").Value = 1.02080289422202
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046181004209715
$ws.Range("D7").Value = 1.050572666951689
$ws.Range("E7").Value = 1.053427871500118
$ws.Range("F7").Value = 1.063352315775545
$ws.Range("I7").Value = 1.039402465664914
$ws.Range("J7").Value = 1.050662355909298
$ws.Range("K7").Value = 1.053019648507139
$ws.Range("L7").Value = 1.05586786892957
$ws.Range("M7").Value = 1.065768350441171
$ws.Range("N7").Value = 1.020716083766723
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044636118832797
$ws.Range("D8").Value = 1.049393142810046
$ws.Range("E8").Value = 1.052072348879072
$ws.Range("F8").Value = 1.061924711198195
$ws.Range("I8").Value = 1.039076343109316
$ws.Range("J8").Value = 1.049582764234902
$ws.Range("K8").Value = 1.052087751892233
$ws.Range("L8").Value = 1.054759696758055
$ws.Range("M8").Value = 1.064585697800971
$ws.Range("N8").Value = 1.020352530354477
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.041912534768207
$ws.Range("D9").Value = 1.047311064500521
$ws.Range("E9").Value = 1.049685049290053
$ws.Range("F9").Value = 1.059409856885344
$ws.Range("I9").Value = 1.038486872757717
$ws.Range("J9").Value = 1.047675060771918
$ws.Range("K9").Value = 1.050437180310324
$ws.Range("L9").Value = 1.052803621553143
$ws.Range("M9").Value = 1.062497913293325
$ws.Range("N9").Value = 1.019709258847602
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040095858710557
$ws.Range("D10").Value = 1.045920577227184
$ws.Range("E10").Value = 1.048094366823236
$ws.Range("F10").Value = 1.057733776350814
$ws.Range("I10").Value = 1.03808394184481
$ws.Range("J10").Value = 1.046399646120801
$ws.Range("K10").Value = 1.049331101218261
$ws.Range("L10").Value = 1.051497302980386
$ws.Range("M10").Value = 1.061103478519075
$ws.Range("N10").Value = 1.019278627010793
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039308973027242
$ws.Range("D11").Value = 1.045317899852219
$ws.Range("E11").Value = 1.047405778587477
$ws.Range("F11").Value = 1.057008123859343
$ws.Range("I11").Value = 1.03790710759229
$ws.Range("J11").Value = 1.045846512639425
$ws.Range("K11").Value = 1.048850800670264
$ws.Range("L11").Value = 1.050931107720998
$ws.Range("M11").Value = 1.060499054956047
$ws.Range("N11").Value = 1.019091733693821
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039016648914366
$ws.Range("D12").Value = 1.045093949995188
$ws.Range("E12").Value = 1.047150033759031
$ws.Range("F12").Value = 1.056738598720725
$ws.Range("I12").Value = 1.037841068399889
$ws.Range("J12").Value = 1.045640922647226
$ws.Range("K12").Value = 1.048672190674376
$ws.Range("L12").Value = 1.050720714158919
$ws.Range("M12").Value = 1.060274450704316
$ws.Range("N12").Value = 1.019022248862541
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039079355254776
$ws.Range("D13").Value = 1.045141992018052
$ws.Range("E13").Value = 1.047204890659821
$ws.Range("F13").Value = 1.056796412133307
$ws.Range("I13").Value = 1.0378552501011
$ws.Range("J13").Value = 1.045685028366086
$ws.Range("K13").Value = 1.048710512409194
$ws.Range("L13").Value = 1.050765848078956
$ws.Range("M13").Value = 1.060322633362577
$ws.Range("N13").Value = 1.019037156509057
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.039284810245465
$ws.Range("D14").Value = 1.045299389898008
$ws.Range("E14").Value = 1.047384638083274
$ws.Range("F14").Value = 1.056985844535653
$ws.Range("I14").Value = 1.037901656016331
$ws.Range("J14").Value = 1.045829521195578
$ws.Range("K14").Value = 1.048836040896947
$ws.Range("L14").Value = 1.050913718229582
$ws.Range("M14").Value = 1.060480491028711
$ws.Range("N14").Value = 1.019085991369488
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.039411392529132
$ws.Range("D15").Value = 1.04539635616737
$ws.Range("E15").Value = 1.047495390026211
$ws.Range("F15").Value = 1.057102561990951
$ws.Range("I15").Value = 1.037930201176883
$ws.Range("J15").Value = 1.045918530535751
$ws.Range("K15").Value = 1.048913355963008
$ws.Range("L15").Value = 1.051004814833094
$ws.Range("M15").Value = 1.060577739814105
$ws.Range("N15").Value = 1.019116071608232
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04014807624472
$ws.Range("D16").Value = 1.045960562531865
$ws.Range("E16").Value = 1.048140070052401
$ws.Range("F16").Value = 1.057781937591566
$ws.Range("I16").Value = 1.038095627964832
$ws.Range("J16").Value = 1.046436337358992
$ws.Range("K16").Value = 1.049362948465546
$ws.Range("L16").Value = 1.051534867833142
$ws.Range("M16").Value = 1.061143578914685
$ws.Range("N16").Value = 1.019291021493917
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.040610109533096
$ws.Range("D17").Value = 1.046314316769289
$ws.Range("E17").Value = 1.048544510574196
$ws.Range("F17").Value = 1.058208118600256
$ws.Range("I17").Value = 1.038198763212842
$ws.Range("J17").Value = 1.046760910271346
$ws.Range("K17").Value = 1.049644601226944
$ws.Range("L17").Value = 1.051867208149447
$ws.Range("M17").Value = 1.061498347276946
$ws.Range("N17").Value = 1.019400648506613
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.040879581442395
$ws.Range("D18").Value = 1.046520599002683
$ws.Range("E18").Value = 1.048780432052917
$ws.Range("F18").Value = 1.058456712609434
$ws.Range("I18").Value = 1.038258692319826
$ws.Range("J18").Value = 1.046950144221374
$ws.Range("K18").Value = 1.049808753235291
$ws.Range("L18").Value = 1.052061003546034
$ws.Range("M18").Value = 1.061705217488437
$ws.Range("N18").Value = 1.019464550913006
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.040971460332051
$ws.Range("D19").Value = 1.046590926321086
$ws.Range("E19").Value = 1.048860878282935
$ws.Range("F19").Value = 1.058541478433082
$ws.Range("I19").Value = 1.03827908793841
$ws.Range("J19").Value = 1.047014653895006
$ws.Range("K19").Value = 1.049864702575782
$ws.Range("L19").Value = 1.052127073784747
$ws.Range("M19").Value = 1.061775744698413
$ws.Range("N19").Value = 1.019486333004086
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.040560540258624
$ws.Range("D20").Value = 1.046276368147712
$ws.Range("E20").Value = 1.048501116059829
$ws.Range("F20").Value = 1.05816239239894
$ws.Range("I20").Value = 1.03818772135441
$ws.Range("J20").Value = 1.046726095357328
$ws.Range("K20").Value = 1.049614396133835
$ws.Range("L20").Value = 1.051831556674649
$ws.Range("M20").Value = 1.061460290217644
$ws.Range("N20").Value = 1.01938889082597
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039224309982309
$ws.Range("D21").Value = 1.045253042595909
$ws.Range("E21").Value = 1.047331706193697
$ws.Range("F21").Value = 1.056930061020126
$ws.Range("I21").Value = 1.03788800043394
$ws.Range("J21").Value = 1.045786975297663
$ws.Range("K21").Value = 1.048799081548813
$ws.Range("L21").Value = 1.050870176464575
$ws.Range("M21").Value = 1.060434008492591
$ws.Range("N21").Value = 1.019071612497734
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03838393822367
$ws.Range("D22").Value = 1.044609124112225
$ws.Range("E22").Value = 1.046596610557763
$ws.Range("F22").Value = 1.056155328931732
$ws.Range("I22").Value = 1.037697499554692
$ws.Range("J22").Value = 1.04519575150723
$ws.Range("K22").Value = 1.048285274658959
$ws.Range("L22").Value = 1.050265235502829
$ws.Range("M22").Value = 1.059788197771026
$ws.Range("N22").Value = 1.018871754716415
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.038829457644829
$ws.Range("D23").Value = 1.044950526263251
$ws.Range("E23").Value = 1.046986283852695
$ws.Range("F23").Value = 1.056566021396874
$ws.Range("I23").Value = 1.037798682459764
$ws.Range("J23").Value = 1.045509242898573
$ws.Range("K23").Value = 1.048557766028746
$ws.Range("L23").Value = 1.050585972176144
$ws.Range("M23").Value = 1.060130606295264
$ws.Range("N23").Value = 1.01897773846406
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.040582938554625
$ws.Range("D24").Value = 1.046293515671052
$ws.Range("E24").Value = 1.048520724116558
$ws.Range("F24").Value = 1.05818305407035
$ws.Range("I24").Value = 1.038192711399112
$ws.Range("J24").Value = 1.046741826977436
$ws.Range("K24").Value = 1.04962804492064
$ws.Range("L24").Value = 1.051847666205262
$ws.Range("M24").Value = 1.061477486750231
$ws.Range("N24").Value = 1.019394203742885
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042616808666903
$ws.Range("D25").Value = 1.047849760037479
$ws.Range("E25").Value = 1.050302071568302
$ws.Range("F25").Value = 1.060059918811322
$ws.Range("I25").Value = 1.038641018861474
$ws.Range("J25").Value = 1.048168882164471
$ws.Range("K25").Value = 1.050864895720954
$ws.Range("L25").Value = 1.053309711727438
$ws.Range("M25").Value = 1.063038108317687
$ws.Range("N25").Value = 1.01987587393532
